$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

# The "CasesTab" (row 2) and "SamplesTab" (row 3) rows still pointed at the
# retired Lymphoma / BCellLymphoma test-data workbooks; point their Neo4j /
# Web data-file columns at the OsteoSarcoma files used by the rest of the
# sheet (rows 4 and 5 already reference them).
$ws.Range("D2").Value = "TC05_Canine_Filter_Diagnosis-OsteoSarcoma_Neo4jData.xlsx"
$ws.Range("E2").Value = "TC05_Canine_Filter_Diagnosis-OsteoSarcoma_WebData.xlsx"
$ws.Range("D3").Value = "TC05_Canine_Filter_Diagnosis-OsteoSarcoma_Neo4jData.xlsx"
$ws.Range("E3").Value = "TC05_Canine_Filter_Diagnosis-OsteoSarcoma_WebData.xlsx"

# Leave the saved view on the FilesTab row's Web data cell, matching the
# author's last selection when the workbook was saved.
$ws.Range("E4").Select()
